# Fix a small bug in the Newton-Raphson column E calculation.
# Column E held a divergent/incorrect iteration; recomputing it correctly
# makes it converge (quadratically) to the root ~1.560040682404455.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    1.546978102189781,
    1.559853335846044,
    1.560040640984931,
    1.560040682404453,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455,
    1.560040682404455
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $values[$i]
}
